$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
# Columns B-S (2-19) keep the original bold/bordered header style (s=1);
# only the text shifts (COD_in is newly inserted after biorefinery, the rest slide
# right by one and a couple are renamed: CAPEX_frac_reduction -> CAPEX_reduction,
# electricity_WWT_frac_reduction -> electricity_WWT_frac_new, etc.)
$ws.Cells.Item(1, 2).Value = "COD_in"
$ws.Cells.Item(1, 3).Value = "MPSP_exist"
$ws.Cells.Item(1, 4).Value = "MPSP_new"
$ws.Cells.Item(1, 5).Value = "MPSP_RIN"
$ws.Cells.Item(1, 6).Value = "MPSP_no_WWT"
$ws.Cells.Item(1, 7).Value = "MPSP_new_frac_reduction"
$ws.Cells.Item(1, 8).Value = "MPSP_RIN_frac_reduction"
$ws.Cells.Item(1, 9).Value = "GWP_exist"
$ws.Cells.Item(1, 10).Value = "GWP_new"
$ws.Cells.Item(1, 11).Value = "GWP_RIN"
$ws.Cells.Item(1, 12).Value = "GWP_no_WWT"
$ws.Cells.Item(1, 13).Value = "GWP_new_frac_reduction"
$ws.Cells.Item(1, 14).Value = "GWP_RIN_frac_reduction"
$ws.Cells.Item(1, 15).Value = "CAPEX_WWT_exist"
$ws.Cells.Item(1, 16).Value = "CAPEX_WWT_new"
$ws.Cells.Item(1, 17).Value = "CAPEX_reduction"
$ws.Cells.Item(1, 18).Value = "CAPEX_WWT_frac_exist"
$ws.Cells.Item(1, 19).Value = "CAPEX_WWT_frac_new"

# Columns T-Z (20-26) are brand new columns beyond the original A:S range, so copy
# the header style (s=1) from A1 onto them first, then set their text.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 20))
$ws.Cells.Item(1, 20).Value = "electricity_WWT_exist"
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 21))
$ws.Cells.Item(1, 21).Value = "electricity_WWT_new"
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 22))
$ws.Cells.Item(1, 22).Value = "electricity_WWT_reduction"
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 23))
$ws.Cells.Item(1, 23).Value = "electricity_WWT_frac_exist"
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 24))
$ws.Cells.Item(1, 24).Value = "electricity_WWT_frac_new"
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 25))
$ws.Cells.Item(1, 25).Value = "ECR_exist"
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 26))
$ws.Cells.Item(1, 26).Value = "ECR_new"

# --- Data rows 2-8 ---
# Row 2: cn
$ws.Cells.Item(2, 2).Value = 7768.034139774792
$ws.Cells.Item(2, 3).Value = 1.473153783329886
$ws.Cells.Item(2, 4).Value = 1.446199816919027
$ws.Cells.Item(2, 5).Value = 1.422134422165871
$ws.Cells.Item(2, 6).Value = 1.424284281648668
$ws.Cells.Item(2, 7).Value = 0.01829677710220634
$ws.Cells.Item(2, 8).Value = 0.03463274624913374
$ws.Cells.Item(2, 9).Value = 2.923263617816073
$ws.Cells.Item(2, 10).Value = 2.644121171773649
$ws.Cells.Item(2, 11).Value = 2.637566847799949
$ws.Cells.Item(2, 12).Value = 2.697530143108954
$ws.Cells.Item(2, 13).Value = 0.09549000108685637
$ws.Cells.Item(2, 14).Value = 0.09773212661181893
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 5.429913831060537
$ws.Cells.Item(2, 17).Value = "-inf"
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0.1532012310107469
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 2155.120155273794
$ws.Cells.Item(2, 22).Value = "-inf"
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0.1232868912155048
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0.2972392960452052

# Row 3: sc1g
$ws.Cells.Item(3, 2).Value = 10304.6145587681
$ws.Cells.Item(3, 3).Value = 2.050894463831623
$ws.Cells.Item(3, 4).Value = 2.020746742237938
$ws.Cells.Item(3, 5).Value = 1.998403519166454
$ws.Cells.Item(3, 6).Value = 1.971636485228009
$ws.Cells.Item(3, 7).Value = 0.01469979178614638
$ws.Cells.Item(3, 8).Value = 0.02559417151436513
$ws.Cells.Item(3, 9).Value = -1.222896471692557
$ws.Cells.Item(3, 10).Value = -1.770089521239921
$ws.Cells.Item(3, 11).Value = -1.689088539678057
$ws.Cells.Item(3, 12).Value = -1.594675566921566
$ws.Cells.Item(3, 13).Value = -0.4474565608894246
$ws.Cells.Item(3, 14).Value = -0.3812195707296989
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 8.32989163355418
$ws.Cells.Item(3, 17).Value = "-inf"
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0.0630929066395696
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 2300.780725416786
$ws.Cells.Item(3, 22).Value = "-inf"
$ws.Cells.Item(3, 23).Value = 0
$ws.Cells.Item(3, 24).Value = 0.0563776575988774
$ws.Cells.Item(3, 25).Value = 0
$ws.Cells.Item(3, 26).Value = 0.0945117444224419

# Row 4: oc1g
$ws.Cells.Item(4, 2).Value = 5448.37141229089
$ws.Cells.Item(4, 3).Value = 1.615587648145031
$ws.Cells.Item(4, 4).Value = 1.639234067316844
$ws.Cells.Item(4, 5).Value = 1.612690004648196
$ws.Cells.Item(4, 6).Value = 1.507058860778469
$ws.Cells.Item(4, 7).Value = -0.01463641988038406
$ws.Cells.Item(4, 8).Value = 0.001793553881252956
$ws.Cells.Item(4, 9).Value = -12.69439032538546
$ws.Cells.Item(4, 10).Value = -13.39460032534528
$ws.Cells.Item(4, 11).Value = -13.28662667705617
$ws.Cells.Item(4, 12).Value = -13.19121830018599
$ws.Cells.Item(4, 13).Value = -0.05515900976824201
$ws.Cells.Item(4, 14).Value = -0.04665339070962589
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 9.76550149632873
$ws.Cells.Item(4, 17).Value = "-inf"
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0.0563333682928295
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 2941.198169644301
$ws.Cells.Item(4, 22).Value = "-inf"
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0.0639882443001577
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 0.1819824538676855

# Row 5: cs
$ws.Cells.Item(5, 2).Value = 65172.40056933241
$ws.Cells.Item(5, 3).Value = 2.144805580120606
$ws.Cells.Item(5, 4).Value = 1.963595577174609
$ws.Cells.Item(5, 5).Value = 1.465663329985225
$ws.Cells.Item(5, 6).Value = 2.166488683623705
$ws.Cells.Item(5, 7).Value = 0.08448784571690986
$ws.Cells.Item(5, 8).Value = 0.3166451339133459
$ws.Cells.Item(5, 9).Value = 0.8725135578166892
$ws.Cells.Item(5, 10).Value = -0.3089016795485015
$ws.Cells.Item(5, 11).Value = 0.5041507988903997
$ws.Cells.Item(5, 12).Value = 3.247051593819852
$ws.Cells.Item(5, 13).Value = 1.354036538207467
$ws.Cells.Item(5, 14).Value = 0.4221857134782548
$ws.Cells.Item(5, 15).Value = 49.4000340563707
$ws.Cells.Item(5, 16).Value = 35.48155915589382
$ws.Cells.Item(5, 17).Value = 0.2817503098195118
$ws.Cells.Item(5, 18).Value = 0.2389994070512131
$ws.Cells.Item(5, 19).Value = 0.2016031428819216
$ws.Cells.Item(5, 20).Value = 60100.50482482319
$ws.Cells.Item(5, 21).Value = 10756.40215879221
$ws.Cells.Item(5, 22).Value = 0.8210264258154863
$ws.Cells.Item(5, 23).Value = 0.3442801822388573
$ws.Cells.Item(5, 24).Value = 0.0878924349000694
$ws.Cells.Item(5, 25).Value = 0.1209860652806428
$ws.Cells.Item(5, 26).Value = 0.0175973939139772

# Row 6: sc2g
$ws.Cells.Item(6, 2).Value = 28827.51152410485
$ws.Cells.Item(6, 3).Value = 2.366646053112585
$ws.Cells.Item(6, 4).Value = 2.117953082951988
$ws.Cells.Item(6, 5).Value = 1.946277347440352
$ws.Cells.Item(6, 6).Value = 1.772905209185035
$ws.Cells.Item(6, 7).Value = 0.1050824519507336
$ws.Cells.Item(6, 8).Value = 0.1776221269417828
$ws.Cells.Item(6, 9).Value = 1.859356273548464
$ws.Cells.Item(6, 10).Value = 1.077731629980625
$ws.Cells.Item(6, 11).Value = 1.327932955197307
$ws.Cells.Item(6, 12).Value = 1.927945360782972
$ws.Cells.Item(6, 13).Value = 0.4203737899440635
$ws.Cells.Item(6, 14).Value = 0.2858103774469046
$ws.Cells.Item(6, 15).Value = 58.33156753466881
$ws.Cells.Item(6, 16).Value = 27.57649820683957
$ws.Cells.Item(6, 17).Value = 0.5272457200734449
$ws.Cells.Item(6, 18).Value = 0.2562523187099998
$ws.Cells.Item(6, 19).Value = 0.1436428053998922
$ws.Cells.Item(6, 20).Value = 35525.83764552922
$ws.Cells.Item(6, 21).Value = 6357.360585689559
$ws.Cells.Item(6, 22).Value = 0.8210496639341139
$ws.Cells.Item(6, 23).Value = 0.3390493569922831
$ws.Cells.Item(6, 24).Value = 0.08365463221873599
$ws.Cells.Item(6, 25).Value = 0.2425118990670642
$ws.Cells.Item(6, 26).Value = 0.0354198676445758

# Row 7: oc2g
$ws.Cells.Item(7, 2).Value = 45473.2229900557
$ws.Cells.Item(7, 3).Value = 2.010626010451889
$ws.Cells.Item(7, 4).Value = 1.756446531932137
$ws.Cells.Item(7, 5).Value = 1.471022358009893
$ws.Cells.Item(7, 6).Value = 1.039754329106481
$ws.Cells.Item(7, 7).Value = 0.1264180793436693
$ws.Cells.Item(7, 8).Value = 0.2683759434310312
$ws.Cells.Item(7, 9).Value = 1.494059199829377
$ws.Cells.Item(7, 10).Value = 0.6035464730411675
$ws.Cells.Item(7, 11).Value = 0.9895990489986736
$ws.Cells.Item(7, 12).Value = 2.081246758710466
$ws.Cells.Item(7, 13).Value = 0.5960357707980425
$ws.Cells.Item(7, 14).Value = 0.3376440176455612
$ws.Cells.Item(7, 15).Value = 50.4601152536607
$ws.Cells.Item(7, 16).Value = 27.59282872011248
$ws.Cells.Item(7, 17).Value = 0.4531754717284219
$ws.Cells.Item(7, 18).Value = 0.2040279456748265
$ws.Cells.Item(7, 19).Value = 0.1221805305529623
$ws.Cells.Item(7, 20).Value = 27900.92776762232
$ws.Cells.Item(7, 21).Value = 5029.401505451774
$ws.Cells.Item(7, 22).Value = 0.8197407072861516
$ws.Cells.Item(7, 23).Value = 0.2764618413268506
$ws.Cells.Item(7, 24).Value = 0.06433720810622311
$ws.Cells.Item(7, 25).Value = 0.1534023211718589
$ws.Cells.Item(7, 26).Value = 0.0221352786626753

# Row 8: la
$ws.Cells.Item(8, 2).Value = 91190.00268360044
$ws.Cells.Item(8, 3).Value = 1.375168058836285
$ws.Cells.Item(8, 4).Value = 1.302946746890901
$ws.Cells.Item(8, 5).Value = 1.077820610074246
$ws.Cells.Item(8, 6).Value = 1.353326673072742
$ws.Cells.Item(8, 7).Value = 0.05251817149279955
$ws.Cells.Item(8, 8).Value = 0.2162262618386182
$ws.Cells.Item(8, 9).Value = 4.505451498978486
$ws.Cells.Item(8, 10).Value = 4.200286761975906
$ws.Cells.Item(8, 11).Value = 4.480668837191839
$ws.Cells.Item(8, 12).Value = 6.000633948144233
$ws.Cells.Item(8, 13).Value = 0.06773233205856723
$ws.Cells.Item(8, 14).Value = 0.0055005945113971
$ws.Cells.Item(8, 15).Value = 42.13425993396559
$ws.Cells.Item(8, 16).Value = 42.59115229752614
$ws.Cells.Item(8, 17).Value = -0.0108437258486706
$ws.Cells.Item(8, 18).Value = 0.1520458131793466
$ws.Cells.Item(8, 19).Value = 0.1534418563043897
$ws.Cells.Item(8, 20).Value = 76596.63085865349
$ws.Cells.Item(8, 21).Value = 9666.512176416243
$ws.Cells.Item(8, 22).Value = 0.8737997733313596
$ws.Cells.Item(8, 23).Value = 0.2541881249357208
$ws.Cells.Item(8, 24).Value = 0.0412381091990413
$ws.Cells.Item(8, 25).Value = 0.0581715149130677
$ws.Cells.Item(8, 26).Value = 0.0067108247407394

